$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.659.84"
$ws.Range("E2").Value = "  +4.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.635.32"
$ws.Range("E3").Value = "  +4.12%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "202.36"
$ws.Range("E5").Value = "  +8.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "583.46"
$ws.Range("E6").Value = "  +2.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.628.70"
$ws.Range("E7").Value = "  +3.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.625"
$ws.Range("E8").Value = "  +4.01%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.689"
$ws.Range("E10").Value = "  +5.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "61.12"
$ws.Range("E11").Value = "  +17.88%  "
$ws.Range("E12").Value = "  +6.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000289"
$ws.Range("E13").Value = "  +14.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.14"
$ws.Range("E14").Value = "  +6.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.200.78"
$ws.Range("E15").Value = "  +3.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.629.10"
$ws.Range("E16").Value = "  +3.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.22"
$ws.Range("E18").Value = "  +7.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.55"
$ws.Range("E19").Value = "  +6.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.428.49"
$ws.Range("E20").Value = "  +4.67%  "
$ws.Range("E21").Value = "  +4.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "407.87"
$ws.Range("E22").Value = "  +6.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.96"
$ws.Range("E23").Value = "  +20.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.30"
$ws.Range("E24").Value = "  +2.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.13"
$ws.Range("E25").Value = "  +2.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.95"
$ws.Range("E26").Value = "  +4.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.00"
$ws.Range("E27").Value = "  +17.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.78"
$ws.Range("E28").Value = "  +6.08%  "
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.51"
$ws.Range("E30").Value = "  +10.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.91"
$ws.Range("E31").Value = "  +13.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.95"
$ws.Range("E32").Value = "  +5.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "682.82"
$ws.Range("E33").Value = "  +12.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.37"
$ws.Range("E34").Value = "  +4.03%  "
$ws.Range("E35").Value = "  +4.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.89"
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.17"
$ws.Range("E37").Value = "  +3.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.422"
$ws.Range("E38").Value = "  +8.99%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0781"
$ws.Range("E40").Value = "  +6.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.24"
$ws.Range("E41").Value = "  +18.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.201.60"
$ws.Range("E42").Value = "  +9.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  +5.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.73"
$ws.Range("E44").Value = "  +12.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.996"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("E46").Value = "  +28.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.87"
$ws.Range("E47").Value = "  +16.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0421"
$ws.Range("E48").Value = "  +6.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.133"
$ws.Range("E49").Value = "  +3.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.87"
$ws.Range("E50").Value = "  +7.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.11"
$ws.Range("E51").Value = "  +1.98%  "
